$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM line: "profile 20 x 40" (lidar mount support profile), row 29
$ws.Range("B29").Value = "profile 20 x 40"
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 116.84

# H29/D29/E29/I29/J29 formulas already exist (shared formulas spanning
# rows 28:32) and will recalculate automatically once F29/G29/C29 are filled.

# Link cell + hyperlink for the new row
$ws.Range("M29").Value = "link"
$ws.Hyperlinks.Add($ws.Range("M29"), "https://www.motedis.co.uk/shop/Aluminium-Profile-20x40-I-type-slot-6", "", "", "link")
$ws.Range("M29").Style = $ws.Range("M28").Style

# Move the active selection as recorded after the edit session
$ws.Range("M20").Select()
